$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '22.395.57'
Set-TextValue "E2" '  -4.41%  '
Set-TextValue "D3" '1.569.73'
Set-TextValue "E3" '  -4.63%  '
Set-TextValue "E4" '  -0.09%  '
Set-TextValue "E5" '  -0.07%  '
Set-TextValue "D6" '290.90'
Set-TextValue "E6" '  -2.55%  '
Set-TextValue "D7" '0.3655'
Set-TextValue "E7" '  -3.42%  '
Set-TextValue "D8" '49.34'
Set-TextValue "E8" '  -1.02%  '
Set-TextValue "D9" '0.3374'
Set-TextValue "E9" '  -5.05%  '
Set-TextValue "D10" '1.169'
Set-TextValue "E10" '  -4.02%  '
Set-TextValue "D11" '0.07580'
Set-TextValue "E11" '  -6.28%  '
Set-TextValue "E12" '  -0.10%  '
Set-TextValue "D13" '21.11'
Set-TextValue "E13" '  -4.08%  '
Set-TextValue "D14" '6.051'
Set-TextValue "E14" '  -5.28%  '
Set-TextValue "D15" '6.869'
Set-TextValue "E15" '  -6.47%  '
Set-TextValue "E16" '  -4.75%  '
Set-TextValue "D17" '1.570.20'
Set-TextValue "E17" '  -4.82%  '
Set-TextValue "D18" '89.03'
Set-TextValue "E18" '  -8.70%  '
Set-TextValue "D19" '0.06726'
Set-TextValue "E19" '  -3.24%  '
Set-TextValue "E20" '  -0.02%  '
Set-TextValue "D21" '6.256'
Set-TextValue "D22" '16.43'
Set-TextValue "E22" '  -5.09%  '
Set-TextValue "D23" '0.5246'
Set-TextValue "E23" '  -9.03%  '
Set-TextValue "D24" '11.98'
Set-TextValue "E24" '  -3.30%  '
Set-TextValue "D25" '22.415.20'
Set-TextValue "E25" '  -4.41%  '
Set-TextValue "D26" '2.386'
Set-TextValue "E26" '  -4.29%  '
Set-TextValue "D27" '2.987'
Set-TextValue "E27" '  +2.95%  '
Set-TextValue "D28" '19.87'
Set-TextValue "E28" '  -5.02%  '
Set-TextValue "D29" '144.60'
Set-TextValue "E29" '  -5.26%  '
Set-TextValue "D30" '4.980'
Set-TextValue "E30" '  -4.53%  '
Set-TextValue "D31" '125.08'
Set-TextValue "E31" '  -5.81%  '
Set-TextValue "D32" '1.746.19'
Set-TextValue "B33" 'ImmutableX'
Set-TextValue "C33" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D33" '1.037'
Set-TextValue "E33" '  +3.53%  '
Set-TextValue "B34" 'Filecoin'
Set-TextValue "C34" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D34" '6.290'
Set-TextValue "E34" '  -8.92%  '
Set-TextValue "D35" '1.970'
Set-TextValue "E35" '  -7.02%  '
Set-TextValue "E36" '  -11.20%  '
Set-TextValue "D37" '0.08423'
Set-TextValue "E37" '  -3.40%  '
Set-TextValue "D38" '0.02547'
Set-TextValue "E38" '  -6.05%  '
Set-TextValue "D39" '0.2304'
Set-TextValue "E39" '  -4.87%  '
Set-TextValue "D40" '0.06538'
Set-TextValue "E40" '  -3.41%  '
Set-TextValue "D41" '5.496'
Set-TextValue "E41" '  -6.90%  '
Set-TextValue "E42" '  -8.81%  '
Set-TextValue "D43" '1.250'
Set-TextValue "E43" '  -3.89%  '
Set-TextValue "D44" '0.6385'
Set-TextValue "E44" '  -6.98%  '
Set-TextValue "D45" '14.55'
Set-TextValue "E45" '  -6.12%  '
Set-TextValue "E46" '  -0.11%  '
Set-TextValue "D47" '0.6017'
Set-TextValue "E47" '  -5.33%  '
Set-TextValue "D48" '3.769'
Set-TextValue "E48" '  -3.63%  '
Set-TextValue "D49" '2.129'
Set-TextValue "E49" '  -5.22%  '
Set-TextValue "B50" 'EOS'
Set-TextValue "C50" 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue "D50" '1.218'
Set-TextValue "E50" '  +3.58%  '
Set-TextValue "B51" 'Quant'
Set-TextValue "C51" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D51" '122.01'
Set-TextValue "E51" '  -3.99%  '
